$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header label for column L ("Q5" quiz column)
$ws.Range("L1").Value = "Q5"

# New Q5 scores per student (rows 2-16)
$ws.Range("L2").Value = 7.42
$ws.Range("L3").Value = 8.85
$ws.Range("L4").Value = 8
$ws.Range("L5").Value = 6
$ws.Range("L6").Value = 9.14
$ws.Range("L7").Value = 9.42
$ws.Range("L8").Value = 0
$ws.Range("L9").Value = 9.14
$ws.Range("L10").Value = 9.42
$ws.Range("L11").Value = 7.42
$ws.Range("L12").Value = 9.42
$ws.Range("L13").Value = 0
$ws.Range("L14").Value = 6.85
$ws.Range("L15").Value = 7.71
$ws.Range("L16").Value = 8.57

# Match the author's last UI selection
$ws.Range("B34").Select()
